$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 911.2727
$ws.Range("I28").Value = 597.1
$ws.Range("J28").Value = 4053
$ws.Range("K28").Value = 597.1
$ws.Range("L28").Value = 4053
$ws.Range("M28").Value = -112.1
$ws.Range("N28").Value = -5023
$ws.Range("H87").Value = 54999.5
$ws.Range("J87").Value = 54999.5
$ws.Range("L87").Value = 54999.5
$ws.Range("N87").Value = -57495.5
$ws.Range("H90").Value = 54999.5
$ws.Range("J90").Value = 54999.5
$ws.Range("L90").Value = 164998.5
$ws.Range("N90").Value = -177478.5
$ws.Range("H98").Value = 2562.1052
$ws.Range("I98").Value = 2562.1052
$ws.Range("K98").Value = 2562.1052
$ws.Range("M98").Value = -1064.1052
$ws.Range("H101").Value = 400
$ws.Range("J101").Value = 402.66666
$ws.Range("L101").Value = 1207.99998
$ws.Range("N101").Value = -4451.999980000001
$ws.Range("H103").Value = 809.1429000000001
$ws.Range("I103").Value = 485.5
$ws.Range("J103").Value = 1240.6666
$ws.Range("K103").Value = 1456.5
$ws.Range("L103").Value = 3721.9998
$ws.Range("M103").Value = -870.5
$ws.Range("N103").Value = -4893.9998
$ws.Range("H122").Value = 2562.1052
$ws.Range("I122").Value = 2562.1052
$ws.Range("K122").Value = 7686.3156
$ws.Range("M122").Value = -5236.3156
$ws.Range("H132").Value = 1026.4318
$ws.Range("I132").Value = 1035.5714
$ws.Range("K132").Value = 3106.7142
$ws.Range("M132").Value = -576.7142000000003
$ws.Range("H134").Value = 45000
$ws.Range("J134").Value = 45000
$ws.Range("L134").Value = 45000
$ws.Range("N134").Value = -55140
$ws.Range("H138").Value = 2333.2122
$ws.Range("I138").Value = 1738.5238
$ws.Range("J138").Value = 2771.4036
$ws.Range("K138").Value = 5215.5714
$ws.Range("L138").Value = 8314.210800000001
$ws.Range("M138").Value = -75.57139999999981
$ws.Range("N138").Value = -18594.2108

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2331042
$ws.Range("I32").Value = 2473388
$ws.Range("J32").Value = 25034.2
$ws.Range("K32").Value = 2473388
$ws.Range("L32").Value = 25034.2
$ws.Range("M32").Value = -2473101
$ws.Range("N32").Value = -25608.2
$ws.Range("I61").Value = 1305.919
$ws.Range("K61").Value = 1305.919
$ws.Range("M61").Value = -1093.919
$ws.Range("H63").Value = 2459.7778
$ws.Range("J63").Value = 2396.6667
$ws.Range("L63").Value = 2396.6667
$ws.Range("N63").Value = -3768.6667
$ws.Range("H66").Value = 2459.7778
$ws.Range("J66").Value = 2396.6667
$ws.Range("L66").Value = 11983.3335
$ws.Range("N66").Value = -18847.3335
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
$ws.Range("H132").Value = 4248.7075
$ws.Range("I132").Value = 2500.6206
$ws.Range("J132").Value = 8473.25
$ws.Range("K132").Value = 7501.861800000001
$ws.Range("L132").Value = 25419.75
$ws.Range("M132").Value = -4971.861800000001
$ws.Range("N132").Value = -30479.75
$ws.Range("I136").Value = 1305.919
$ws.Range("K136").Value = 3917.757000000001
$ws.Range("M136").Value = -1367.757000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2199.5
$ws.Range("J94").Value = 4320.231
$ws.Range("L94").Value = 4320.231
$ws.Range("N94").Value = -5222.231
$ws.Range("H107").Value = 102274220
$ws.Range("I107").Value = 102274220
$ws.Range("K107").Value = 102274220
$ws.Range("M107").Value = -102272300

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5873.0386
$ws.Range("I58").Value = 3062.44
$ws.Range("J58").Value = 8475.444
$ws.Range("K58").Value = 3062.44
$ws.Range("L58").Value = 8475.444
$ws.Range("M58").Value = -2859.44
$ws.Range("N58").Value = -8881.444
$ws.Range("H105").Value = 6495130.5
$ws.Range("I105").Value = 10205005
$ws.Range("K105").Value = 10205005
$ws.Range("M105").Value = -10203258
$ws.Range("H132").Value = 3428.6135
$ws.Range("J132").Value = 7748.9
$ws.Range("L132").Value = 23246.7
$ws.Range("N132").Value = -28306.7
$ws.Range("H134").Value = 3524.5862
$ws.Range("I134").Value = 1885.4736
$ws.Range("J134").Value = 6638.9
$ws.Range("K134").Value = 5656.4208
$ws.Range("L134").Value = 19916.7
$ws.Range("M134").Value = -3121.4208
$ws.Range("N134").Value = -24986.7
$ws.Range("H136").Value = 5873.0386
$ws.Range("I136").Value = 3062.44
$ws.Range("J136").Value = 8475.444
$ws.Range("K136").Value = 9187.32
$ws.Range("L136").Value = 25426.332
$ws.Range("M136").Value = -6637.32
$ws.Range("N136").Value = -30526.332

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 41670916
$ws.Range("I80").Value = 26319684
$ws.Range("J80").Value = 100005600
$ws.Range("K80").Value = 78959052
$ws.Range("L80").Value = 300016800
$ws.Range("M80").Value = -78958116
$ws.Range("N80").Value = -300018672
$ws.Range("H83").Value = 41670916
$ws.Range("I83").Value = 26319684
$ws.Range("J83").Value = 100005600
$ws.Range("K83").Value = 236877156
$ws.Range("L83").Value = 900050400
$ws.Range("M83").Value = -236872476
$ws.Range("N83").Value = -900059760
$ws.Range("H113").Value = 7459.933
$ws.Range("I113").Value = 1891.25
$ws.Range("J113").Value = 9484.909
$ws.Range("K113").Value = 5673.75
$ws.Range("L113").Value = 28454.727
$ws.Range("M113").Value = -3503.75
$ws.Range("N113").Value = -32794.727
$ws.Range("H122").Value = 1886828.5
$ws.Range("J122").Value = 856.7143
$ws.Range("L122").Value = 7710.428699999999
$ws.Range("N122").Value = -12610.4287
$ws.Range("H129").Value = 9923202
$ws.Range("J129").Value = 12976144
$ws.Range("L129").Value = 38928432
$ws.Range("N129").Value = -38938432
$ws.Range("H137").Value = 87273.5
$ws.Range("I137").Value = 103202.5
$ws.Range("J137").Value = 75895.64
$ws.Range("K137").Value = 309607.5
$ws.Range("L137").Value = 227686.92
$ws.Range("M137").Value = -304507.5
$ws.Range("N137").Value = -237886.92

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2449.7646
$ws.Range("I80").Value = 2504.9
$ws.Range("J80").Value = 2371
$ws.Range("K80").Value = 2504.9
$ws.Range("L80").Value = 2371
$ws.Range("M80").Value = -1506.9
$ws.Range("N80").Value = -4367
$ws.Range("H83").Value = 2449.7646
$ws.Range("I83").Value = 2504.9
$ws.Range("J83").Value = 2371
$ws.Range("K83").Value = 12524.5
$ws.Range("L83").Value = 11855
$ws.Range("M83").Value = -7532.5
$ws.Range("N83").Value = -21839
$ws.Range("H122").Value = 2338040.5
$ws.Range("I122").Value = 2898478.2
$ws.Range("K122").Value = 8695434.600000001
$ws.Range("M122").Value = -8692984.600000001
$ws.Range("H132").Value = 4473.067
$ws.Range("I132").Value = 5150.25
$ws.Range("K132").Value = 15450.75
$ws.Range("M132").Value = -12920.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5561075.5
$ws.Range("I61").Value = 10003407
$ws.Range("J61").Value = 8161.5
$ws.Range("K61").Value = 10003407
$ws.Range("L61").Value = 8161.5
$ws.Range("M61").Value = -10003205
$ws.Range("N61").Value = -8565.5
$ws.Range("H68").Value = 6041.5835
$ws.Range("I68").Value = 3357.1428
$ws.Range("J68").Value = 9799.799999999999
$ws.Range("K68").Value = 3357.1428
$ws.Range("L68").Value = 9799.799999999999
$ws.Range("M68").Value = -2608.1428
$ws.Range("N68").Value = -11297.8
$ws.Range("H71").Value = 6041.5835
$ws.Range("I71").Value = 3357.1428
$ws.Range("J71").Value = 9799.799999999999
$ws.Range("K71").Value = 16785.714
$ws.Range("L71").Value = 48999
$ws.Range("M71").Value = -13041.714
$ws.Range("N71").Value = -56487
$ws.Range("H113").Value = 5561075.5
$ws.Range("I113").Value = 10003407
$ws.Range("J113").Value = 8161.5
$ws.Range("K113").Value = 10003407
$ws.Range("L113").Value = 8161.5
$ws.Range("M113").Value = -10001237
$ws.Range("N113").Value = -12501.5
$ws.Range("H122").Value = 2873.282
$ws.Range("I122").Value = 2117.6155
$ws.Range("J122").Value = 4384.615
$ws.Range("K122").Value = 6352.8465
$ws.Range("L122").Value = 13153.845
$ws.Range("M122").Value = -3902.8465
$ws.Range("N122").Value = -18053.845
$ws.Range("H136").Value = 10077.933
$ws.Range("I136").Value = 3200.6453
$ws.Range("J136").Value = 17692.072
$ws.Range("K136").Value = 9601.9359
$ws.Range("L136").Value = 53076.216
$ws.Range("M136").Value = -7051.9359
$ws.Range("N136").Value = -58176.216

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 24.363636
$ws.Range("I21").Value = 24.363636
$ws.Range("K21").Value = 24.363636
$ws.Range("M21").Value = 210.636364
$ws.Range("H26").Value = 10012
$ws.Range("I26").Value = 10012
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 10012
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -9719
$ws.Range("N26").ClearContents()
$ws.Range("H28").Value = 5000
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H35").Value = 24.363636
$ws.Range("I35").Value = 24.363636
$ws.Range("K35").Value = 24.363636
$ws.Range("M35").Value = 265.636364
$ws.Range("H122").Value = 110658.516
$ws.Range("I122").Value = 175609.39
$ws.Range("J122").Value = 3953.5
$ws.Range("K122").Value = 526828.17
$ws.Range("L122").Value = 11860.5
$ws.Range("M122").Value = -524378.17
$ws.Range("N122").Value = -16760.5
$ws.Range("H132").Value = 5344.081
$ws.Range("I132").Value = 5942.476
$ws.Range("K132").Value = 17827.428
$ws.Range("M132").Value = -15297.428
$ws.Range("H135").Value = 73999.5
$ws.Range("J135").Value = 73999.5
$ws.Range("L135").Value = 73999.5
$ws.Range("N135").Value = -84139.5
